$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $targetId) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.Id -eq $targetId) {
            return $cand
        }
    }
    return $null
}

# --- Shape id=88 "Add-Ons" background roundRect: reposition/resize only ---
$sh88 = Get-ShapeById $s 88
$sh88.Left = 15.961742283464567
$sh88.Top = 358.5655118110236
$sh88.Width = 444.0826971653543
$sh88.Height = 185.89023622047245

# --- Shape id=93 "+Cheese / +Jalapeños / +Chili / +Butter": drop the "+Butter" line, reposition/resize ---
$sh93 = Get-ShapeById $s 93
$tr93 = $sh93.TextFrame.TextRange
$tr93.Text = "+Cheese`r+Jalapeños`r+Chili"
$sh93.Left = 14.803002125984252
$sh93.Top = 413.7387601574803
$sh93.Width = 318.78520685039365
$sh93.Height = 119.95669291338582

# --- Shape id=6 (four "$0.50" price lines): drop the first "$0.50" line, reposition/resize ---
$sh6 = Get-ShapeById $s 6
$tr6 = $sh6.TextFrame.TextRange
$tr6.Paragraphs(1).Delete()
$sh6.Left = 344.164114488189
$sh6.Top = 412.775905511811
$sh6.Width = 102.63771653543307
$sh6.Height = 119.95669291338582
